$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ASSESSMENTS")
$ws2 = $wb.Worksheets.Item("INDUSTRIES")

# Rename "industry" references to "site" on the INDUSTRIES sheet
$ws2.Range("B3").Value = "Site "
$ws2.Range("A1").Value = "SITE"
$ws2.Range("E1").Value = "SUB-SUPPLIERS"

# Update the active sheet / selections to match the new state
$ws1.Activate() | Out-Null
$ws1.Range("B12").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("E1:J1").Select() | Out-Null
